# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (B16:G33) is re-organized: instead of
# grouping all periods (2305-2308) under each worker, it is now grouped
# by period (2305, 2306, 2307, 2308) with all five workers listed under
# each period. The "Valor Mora" (G) totals for CRISTIAN CAMILO CASANOVA
# CARRASCAL and AMINTA ROSA CASTRO MERCADO are also corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Tipo Doc, N Doc, Nombre, Periodo Mora, Valor Mora, Salario Basico
$rows = @(
    @("CC", "1143358292", "ALDAIR MAURICIO MARTINEZ MONTALVO",  "2305", 46400,  1160000),
    @("CC", "1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL",  "2305", 160000, 4347826),
    @("CC", "1128050143", "LEIDYS DEL CARMEN SARA PAJARO",       "2305", 80000,  2000000),
    @("CC", "1237439154", "YESSIKA MARIA LOPEZ YI",              "2305", 46400,  1160000),
    @("CC", "1002428397", "AMINTA ROSA CASTRO MERCADO",          "2305", 50160,  1160000),

    @("CC", "1143358292", "ALDAIR MAURICIO MARTINEZ MONTALVO",  "2306", 46400,  1160000),
    @("CC", "1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL",  "2306", 160000, 4347826),
    @("CC", "1128050143", "LEIDYS DEL CARMEN SARA PAJARO",       "2306", 80000,  2000000),
    @("CC", "1237439154", "YESSIKA MARIA LOPEZ YI",              "2306", 46400,  1160000),
    @("CC", "1002428397", "AMINTA ROSA CASTRO MERCADO",          "2306", 50160,  1160000),

    @("CC", "1143358292", "ALDAIR MAURICIO MARTINEZ MONTALVO",  "2307", 36341,  1160000),
    @("CC", "1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL",  "2307", 160000, 4347826),
    @("CC", "1128050143", "LEIDYS DEL CARMEN SARA PAJARO",       "2307", 80000,  2000000),
    @("CC", "1237439154", "YESSIKA MARIA LOPEZ YI",              "2307", 46400,  1160000),
    @("CC", "1002428397", "AMINTA ROSA CASTRO MERCADO",          "2307", 50160,  1160000),

    @("CC", "1143332822", "CRISTIAN CAMILO CASANOVA CARRASCAL",  "2308", 160000, 4347826),
    @("CC", "1128050143", "LEIDYS DEL CARMEN SARA PAJARO",       "2308", 80000,  2000000),
    @("CC", "1237439154", "YESSIKA MARIA LOPEZ YI",              "2308", 46400,  1160000)
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]   # B Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $row[1]   # C N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[2]   # D Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[3]   # E Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[4]   # F Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[5]   # G Salario Basico
    $r++
}
